$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "org.openmrs.Patient,org.openmrs.module.patientflags.PatientFlag,org.openmrs.module.drools.calculation.Operator, static org.openmrs.module.drools.utils.DroolsDateUtils.daysAgo"

$ws.Range("H11").Value = "not PatientFlag(patient == `$patient)"
$ws.Range("I11").Value = "`$flag: PatientFlag(patient == `$patient, message == `$param)"
$ws.Range("J11").Value = "insert(new PatientFlag(`$patient, null, `$param));"

$ws.Range("C2").Select()
